# Updates the "properties correspondence table for biotic elements":
# applies an extensiveness penalty in the similarity measure by bumping
# every correction value of 5 or 7 in the matrix up to 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "G1","L1","H2","M2","I3","N3","J4","O4","K5","P5",
    "B6","L6","C7","M7","D8","N8","E9","O9","F10","P10",
    "B11","G11","C12","H12","D13","I13","E14","J14","F15","K15",
    "T16","W16","Y16","U17","X17","Y18","Q19","V19","W19","R20",
    "X20","S21","Y21","Q22","T22","R23","U23","S24","V24"
)

foreach ($cell in $cells) {
    $ws.Range($cell).Value = 8
}

# H2 additionally gets flagged with a red font to highlight the change.
$ws.Range("H2").Font.Color = 255

# Leave column A selected, as when reviewing the updated table.
$ws.Columns("A").Select() | Out-Null
